$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.410411333333334
$ws.Range("H2").Value = 25.231234
$ws.Range("I2").Value = 0.004700308673809217
$ws.Range("J2").Value = 0.004733981166790757
$ws.Range("M2").Value = 2.213484666666667
$ws.Range("N2").Value = 6.640453999999999
$ws.Range("O2").Value = 0.1355468747528199
$ws.Range("P2").Value = 0.1561506731221149
$ws.Range("Q2").Value = 18.61631652669289
$ws.Range("R2").Value = 167.546848740236
$ws.Range("S2").Value = 0.0006371121511084109
$ws.Range("T2").Value = 0.0007392143457417917
$ws.Range("G3").Value = 8.410411333333334
$ws.Range("H3").Value = 25.231234
$ws.Range("I3").Value = 0.004700308673809217
$ws.Range("J3").Value = 0.004733981166790757
$ws.Range("O3").Value = 0.4686085613268789
$ws.Range("P3").Value = 0.5398393907304434
$ws.Range("Q3").Value = 64.35976720737979
$ws.Range("R3").Value = 579.237904866418
$ws.Range("S3").Value = 0.002202604885425987
$ws.Range("T3").Value = 0.002555589508809716
$ws.Range("G4").Value = 8.410411333333334
$ws.Range("H4").Value = 25.231234
$ws.Range("I4").Value = 0.004700308673809217
$ws.Range("J4").Value = 0.004733981166790757
$ws.Range("M4").Value = 6.464154000000001
$ws.Range("N4").Value = 12.928308
$ws.Range("O4").Value = 0.3958445639203012
$ws.Range("P4").Value = 0.3040099361474417
$ws.Range("Q4").Value = 54.36619406201201
$ws.Range("R4").Value = 326.1971643720721
$ws.Range("S4").Value = 0.001860591637274819
$ws.Range("T4").Value = 0.00143917731223925
$ws.Range("I5").Value = 0.972865573481432
$ws.Range("J5").Value = 0.979835075160667
$ws.Range("M5").Value = 2.213484666666667
$ws.Range("N5").Value = 6.640453999999999
$ws.Range("O5").Value = 0.1355468747528199
$ws.Range("P5").Value = 0.1561506731221149
$ws.Range("Q5").Value = 3853.188101192365
$ws.Range("R5").Value = 34678.69291073128
$ws.Range("S5").Value = 0.131868888040018
$ws.Range("T5").Value = 0.1530019065349962
$ws.Range("I6").Value = 0.972865573481432
$ws.Range("J6").Value = 0.979835075160667
$ws.Range("O6").Value = 0.4686085613268789
$ws.Range("P6").Value = 0.5398393907304434
$ws.Range("S6").Value = 0.4558931367535828
$ws.Range("T6").Value = 0.5289535699910527
$ws.Range("I7").Value = 0.972865573481432
$ws.Range("J7").Value = 0.979835075160667
$ws.Range("M7").Value = 6.464154000000001
$ws.Range("N7").Value = 12.928308
$ws.Range("O7").Value = 0.3958445639203012
$ws.Range("P7").Value = 0.3040099361474417
$ws.Range("Q7").Value = 11252.66492791383
$ws.Range("R7").Value = 67515.98956748297
$ws.Range("S7").Value = 0.3851035486878312
$ws.Range("T7").Value = 0.2978795986346181
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.318483333333333
$ws.Range("H8").Value = 3.95545
$ws.Range("I8").Value = 0.0007368579730907598
$ws.Range("J8").Value = 0.0007421367423480953
$ws.Range("M8").Value = 2.213484666666667
$ws.Range("N8").Value = 6.640453999999999
$ws.Range("O8").Value = 0.1355468747528199
$ws.Range("P8").Value = 0.1561506731221149
$ws.Range("Q8").Value = 2.918442641588888
$ws.Range("R8").Value = 26.2659837743
$ws.Range("S8").Value = 0.00009987879538914996
$ws.Range("T8").Value = 0.0001158851518663087
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.318483333333333
$ws.Range("H9").Value = 3.95545
$ws.Range("I9").Value = 0.0007368579730907598
$ws.Range("J9").Value = 0.0007421367423480953
$ws.Range("O9").Value = 0.4686085613268789
$ws.Range("P9").Value = 0.5398393907304434
$ws.Range("Q9").Value = 10.08955175162778
$ws.Range("R9").Value = 90.80596576465
$ws.Range("S9").Value = 0.000345297954672301
$ws.Range("T9").Value = 0.0004006346468278718
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.318483333333333
$ws.Range("H10").Value = 3.95545
$ws.Range("I10").Value = 0.0007368579730907598
$ws.Range("J10").Value = 0.0007421367423480953
$ws.Range("M10").Value = 6.464154000000001
$ws.Range("N10").Value = 12.928308
$ws.Range("O10").Value = 0.3958445639203012
$ws.Range("P10").Value = 0.3040099361474417
$ws.Range("Q10").Value = 8.522879313100001
$ws.Range("R10").Value = 51.1372758786
$ws.Range("S10").Value = 0.0002916812230293088
$ws.Range("T10").Value = 0.0002256169436539149
$ws.Range("G11").Value = 38.182192
$ws.Range("H11").Value = 76.364384
$ws.Range("I11").Value = 0.02133880034277938
$ws.Range("J11").Value = 0.0143277794367718
$ws.Range("M11").Value = 2.213484666666667
$ws.Range("N11").Value = 6.640453999999999
$ws.Range("O11").Value = 0.1355468747528199
$ws.Range("P11").Value = 0.1561506731221149
$ws.Range("Q11").Value = 84.51569653172267
$ws.Range("R11").Value = 507.0941791903359
$ws.Range("S11").Value = 0.002892407697438146
$ws.Range("T11").Value = 0.002237292403397113
$ws.Range("G12").Value = 38.182192
$ws.Range("H12").Value = 76.364384
$ws.Range("I12").Value = 0.02133880034277938
$ws.Range("J12").Value = 0.0143277794367718
$ws.Range("O12").Value = 0.4686085613268789
$ws.Range("P12").Value = 0.5398393907304434
$ws.Range("Q12").Value = 292.1851133306614
$ws.Range("R12").Value = 1753.110679983968
$ws.Range("S12").Value = 0.009999544529071354
$ws.Range("T12").Value = 0.007734699721667063
$ws.Range("G13").Value = 38.182192
$ws.Range("H13").Value = 76.364384
$ws.Range("I13").Value = 0.02133880034277938
$ws.Range("J13").Value = 0.0143277794367718
$ws.Range("M13").Value = 6.464154000000001
$ws.Range("N13").Value = 12.928308
$ws.Range("O13").Value = 0.3958445639203012
$ws.Range("P13").Value = 0.3040099361474417
$ws.Range("Q13").Value = 246.815569145568
$ws.Range("R13").Value = 987.2622765822721
$ws.Range("S13").Value = 0.008446848116269876
$ws.Range("T13").Value = 0.004355787311707622
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6414029999999999
$ws.Range("H14").Value = 1.924209
$ws.Range("I14").Value = 0.0003584595288887479
$ws.Range("J14").Value = 0.0003610274934222114
$ws.Range("M14").Value = 2.213484666666667
$ws.Range("N14").Value = 6.640453999999999
$ws.Range("O14").Value = 0.1355468747528199
$ws.Range("P14").Value = 0.1561506731221149
$ws.Range("Q14").Value = 1.419735705654
$ws.Range("R14").Value = 12.777621350886
$ws.Range("S14").Value = 0.00004858806886623794
$ws.Range("T14").Value = 0.00005637468611346822
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6414029999999999
$ws.Range("H15").Value = 1.924209
$ws.Range("I15").Value = 0.0003584595288887479
$ws.Range("J15").Value = 0.0003610274934222114
$ws.Range("O15").Value = 0.4686085613268789
$ws.Range("P15").Value = 0.5398393907304434
$ws.Range("Q15").Value = 4.908267399777
$ws.Range("R15").Value = 44.174406597993
$ws.Range("S15").Value = 0.0001679772041264669
$ws.Range("T15").Value = 0.0001948968620859858
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6414029999999999
$ws.Range("H16").Value = 1.924209
$ws.Range("I16").Value = 0.0003584595288887479
$ws.Range("J16").Value = 0.0003610274934222114
$ws.Range("M16").Value = 6.464154000000001
$ws.Range("N16").Value = 12.928308
$ws.Range("O16").Value = 0.3958445639203012
$ws.Range("P16").Value = 0.3040099361474417
$ws.Range("Q16").Value = 4.146127768062
$ws.Range("R16").Value = 24.876766608372
$ws.Range("S16").Value = 0.000141894255896043
$ws.Range("T16").Value = 0.0001097559452227574
